$d = $word.ActiveDocument

# --- A. Remove the _GoBack bookmark currently sitting in the first (empty) paragraph ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- B. "12345678 CC" -> "12345678 ID" ---
$d.Content.Find.Execute("12345678 CC", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "12345678 ID", 2)

# --- C. "NIF 111111111" paragraph becomes "VAT" / "/NIF" / <bookmark/> / " 111111111" ---
$pNif = $d.Paragraphs.Item(7)
$pNif.Range.Text = "VAT/NIF 111111111"

# Force a hard run-break between "VAT" and "/NIF" (no formatting difference, so the only
# reliable way to keep two <w:r> elements instead of Word re-merging them is to drop a
# bookmark on the seam and remove it again once the split has taken).
$pNif = $d.Paragraphs.Item(7)
$splitPoint = $pNif.Range.Start + 3
$tmp = $d.Range($splitPoint, $splitPoint)
$d.Bookmarks.Add("TempRunSplit", $tmp)
$d.Bookmarks.Item("TempRunSplit").Delete()

# Re-anchor _GoBack between "/NIF" and " 111111111"
$pNif = $d.Paragraphs.Item(7)
$bmPoint = $pNif.Range.Start + 7
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)

# --- D. Two new paragraphs after the VAT/NIF paragraph ---
$pNif = $d.Paragraphs.Item(7)
$pNif.Range.InsertParagraphAfter()
$pUtr = $d.Paragraphs.Item(8)
$pUtr.Range.Text = "UTR 1234567890"

$pUtr = $d.Paragraphs.Item(8)
$pUtr.Range.InsertParagraphAfter()
$pNino = $d.Paragraphs.Item(9)
$pNino.Range.Text = "NINO DQ123456D"

# --- E. One extra trailing empty paragraph (2 -> 3) ---
$pNino = $d.Paragraphs.Item(9)
$pNino.Range.InsertParagraphAfter()

Write-Host $d.Content.Text
